$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 4 (Ayaz Mohammad Zahir) : fill in newly-scored test data ---
$ws.Range("C4").Value = 13
$ws.Range("D4").Value = 3
$ws.Range("G4").Value = 10
$ws.Range("H4").Value = 2
$ws.Range("K4").Value = 8
$ws.Range("L4").Value = 6
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("S4").Value = 12
$ws.Range("T4").Value = 3

# --- Row 20 : corrected English Correct/Wrong counts ---
$ws.Range("C20").Value = 15
$ws.Range("D20").Value = 0

# --- sheet view / selection cosmetics ---
$ws.Application.ActiveWindow.ScrollColumn = 5
$ws.Range("AF4").Select()
